# Update the cryptos list worksheet with the latest scraped values.
# Values are assigned with a leading apostrophe to force Excel to treat
# numeric-looking strings (prices, percentages) as text, matching the
# original workbook's string-typed cells. Style is reset to Normal
# afterwards so we don't leave a stray "Text" number format behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'72.943.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +3.12%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.982.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.24%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'596.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +11.89%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'161.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +9.63%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.49%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D9").Value = "'0.747"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.64%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +2.32%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'54.54"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.02%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.0000318"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.31%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'10.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +3.57%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'4.622.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.20%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.987.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.21%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +9.65%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'14.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'20.32"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.74%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.34%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'72.614.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.74%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'435.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.24%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +13.30%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'95.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.08%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'3.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -4.27%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'14.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.73%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'4.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +13.74%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'11.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.19%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +1.63%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'10.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.46%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'36.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.09%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'7.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.43%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'13.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +3.43%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.08%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'48.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -4.90%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'667.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.03%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +8.84%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.0₃0899"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +11.19%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.436"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Value = "'1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.12%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'ThetaToken"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'3.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.23%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'WEMIXToken"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'3.34"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +5.21%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'Kaspa"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.145"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.47%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +0.14%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +2.23%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'10.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +6.66%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +0.50%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'3.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +3.63%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.15%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.888.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +9.86%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +1.95%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +4.70%  "
$ws.Range("E51").Style = "Normal"
